# Update the "Förändrad" date column (C) for rows 2-10 from 45174 (2023-09-05)
# to 45175 (2023-09-06) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Value = 45175
